$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-138 all change from serial date 45202 (2023-10-03)
# to serial date 45203 (2023-10-04). Update the value while preserving the
# existing date-formatted style of the cells.
$ws.Range("C2:C138").Value = 45203
